$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "MEC-2B-Mec. Tec. Res. Mat."

# Row 6
$ws.Range("D6").Value = "MEC-2B-Mec. Tec. Res. Mat."
$ws.Range("F6").Value = "-"

# Row 10
$ws.Range("D10").Value = "MEC-2A-Mec. Tec. Res. Mat."

# Row 11
$ws.Range("D11").Value = "MEC-2A-Mec. Tec. Res. Mat."
$ws.Range("F11").Value = "-"

# Row 15
$ws.Range("B15").Value = "-"

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("E18").Value = "MEC-1NB-M.T.R.M."
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("E19").Value = "ELM-1NA-Tecnologias Mecânicas"

# Row 20
$ws.Range("D20").Value = "-"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "MEC-1NA-M.T.R.M."
$ws.Range("C21").Value = "MEC-1NB-M.T.R.M."
$ws.Range("E21").Value = "MEC-1NB-M.T.R.M."
